$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 223.4
$ws.Range("I2").Value = 172.33333
$ws.Range("K2").Value = 172.33333
$ws.Range("M2").Value = -59.33332999999999

$ws.Range("H19").Value = 15253
$ws.Range("I19").Value = 1038.2
$ws.Range("J19").Value = 23150.111
$ws.Range("K19").Value = 1038.2
$ws.Range("L19").Value = 23150.111
$ws.Range("M19").Value = -863.2
$ws.Range("N19").Value = -23500.111

$ws.Range("H100").Value = 2540.3333
$ws.Range("J100").Value = 2377.7778
$ws.Range("L100").Value = 2377.7778
$ws.Range("N100").Value = -3459.7778

$ws.Range("H103").Value = 83740.75
$ws.Range("I103").Value = 100423.9
$ws.Range("J103").Value = 325
$ws.Range("K103").Value = 301271.7
$ws.Range("L103").Value = 975
$ws.Range("M103").Value = -300685.7
$ws.Range("N103").Value = -2147

$ws.Range("H133").Value = 60780
$ws.Range("J133").Value = 60780
$ws.Range("L133").Value = 60780
$ws.Range("N133").Value = -70900

$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws.Range("H138").Value = 2870.6365
$ws.Range("I138").Value = 16598.5
$ws.Range("J138").Value = 2216.9285
$ws.Range("K138").Value = 49795.5
$ws.Range("L138").Value = 6650.7855
$ws.Range("M138").Value = -44655.5
$ws.Range("N138").Value = -16930.7855

$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("M139").ClearContents()
$ws.Range("N139").ClearContents()

$ws.Range("H140").Value = 59690
$ws.Range("J140").Value = 59690
$ws.Range("L140").Value = 59690
$ws.Range("N140").Value = -70050

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 549409.5
$ws.Range("I32").Value = 642122.9
$ws.Range("J32").Value = 16307.75
$ws.Range("K32").Value = 642122.9
$ws.Range("L32").Value = 16307.75
$ws.Range("M32").Value = -641835.9
$ws.Range("N32").Value = -16881.75

$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

$ws.Range("H61").Value = 2727.75
$ws.Range("I61").Value = 1537.3334
$ws.Range("J61").Value = 3124.5557
$ws.Range("K61").Value = 1537.3334
$ws.Range("L61").Value = 3124.5557
$ws.Range("M61").Value = -1325.3334
$ws.Range("N61").Value = -3548.5557

$ws.Range("H74").Value = 1369.1904
$ws.Range("I74").Value = 868.0714
$ws.Range("J74").Value = 2371.4285
$ws.Range("K74").Value = 868.0714
$ws.Range("L74").Value = 2371.4285
$ws.Range("M74").Value = 5.92859999999996
$ws.Range("N74").Value = -4119.4285

$ws.Range("H77").Value = 1369.1904
$ws.Range("I77").Value = 868.0714
$ws.Range("J77").Value = 2371.4285
$ws.Range("K77").Value = 4340.357
$ws.Range("L77").Value = 11857.1425
$ws.Range("M77").Value = 27.64300000000003
$ws.Range("N77").Value = -20593.1425

$ws.Range("H136").Value = 2727.75
$ws.Range("I136").Value = 1537.3334
$ws.Range("J136").Value = 3124.5557
$ws.Range("K136").Value = 4612.0002
$ws.Range("L136").Value = 9373.667099999999
$ws.Range("M136").Value = -2062.0002
$ws.Range("N136").Value = -14473.6671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2942.7715
$ws.Range("I134").Value = 2799.32
$ws.Range("K134").Value = 8397.960000000001
$ws.Range("M134").Value = -5862.960000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 69002
$ws.Range("J4").Value = 69002
$ws.Range("L4").Value = 69002
$ws.Range("N4").Value = -69226

$ws.Range("H58").Value = 1402.2222
$ws.Range("I58").Value = 1140
$ws.Range("J58").Value = 3500
$ws.Range("K58").Value = 1140
$ws.Range("L58").Value = 3500
$ws.Range("M58").Value = -937
$ws.Range("N58").Value = -3906

$ws.Range("H122").Value = 1522.8286
$ws.Range("I122").Value = 1072.1666
$ws.Range("K122").Value = 3216.4998
$ws.Range("M122").Value = -766.4998000000001

$ws.Range("H132").Value = 15153886
$ws.Range("I132").Value = 1012.3333
$ws.Range("J132").Value = 33337334
$ws.Range("K132").Value = 3036.9999
$ws.Range("L132").Value = 100012002
$ws.Range("M132").Value = -506.9998999999998
$ws.Range("N132").Value = -100017062

$ws.Range("H134").Value = 4089.8572
$ws.Range("I134").Value = 1124.8
$ws.Range("J134").Value = 11502.5
$ws.Range("K134").Value = 3374.4
$ws.Range("L134").Value = 34507.5
$ws.Range("M134").Value = -839.3999999999996
$ws.Range("N134").Value = -39577.5

$ws.Range("H136").Value = 1402.2222
$ws.Range("I136").Value = 1140
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 3420
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -870
$ws.Range("N136").Value = -15600

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 20833722
$ws.Range("J34").Value = 27778250
$ws.Range("L34").Value = 83334750
$ws.Range("N34").Value = -83334918

$ws.Range("H131").Value = 675.625
$ws.Range("I131").Value = 394.44446
$ws.Range("J131").Value = 1037.1428
$ws.Range("K131").Value = 1183.33338
$ws.Range("L131").Value = 3111.4284
$ws.Range("M131").Value = 3856.66662
$ws.Range("N131").Value = -13191.4284

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3785.64
$ws.Range("I122").Value = 2404.077
$ws.Range("K122").Value = 7212.231000000001
$ws.Range("M122").Value = -4762.231000000001

$ws.Range("H132").Value = 2770.2222
$ws.Range("I132").Value = 2222
$ws.Range("J132").Value = 4195.6
$ws.Range("K132").Value = 6666
$ws.Range("L132").Value = 12586.8
$ws.Range("M132").Value = -4136
$ws.Range("N132").Value = -17646.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2895.5945
$ws.Range("I132").Value = 2508.9
$ws.Range("J132").Value = 4552.857
$ws.Range("K132").Value = 7526.700000000001
$ws.Range("L132").Value = 13658.571
$ws.Range("M132").Value = -4996.700000000001
$ws.Range("N132").Value = -18718.571

$ws.Range("H136").Value = 2778596.8
$ws.Range("I136").Value = 850.7838
$ws.Range("J136").Value = 7247144.5
$ws.Range("K136").Value = 2552.3514
$ws.Range("L136").Value = 21741433.5
$ws.Range("M136").Value = -2.351400000000012
$ws.Range("N136").Value = -21746533.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 40715
$ws.Range("J133").Value = 40715
$ws.Range("L133").Value = 40715
$ws.Range("N133").Value = -50835

$ws.Range("H136").Value = 3572.9443
$ws.Range("I136").Value = 3385.08
$ws.Range("J136").Value = 3999.9092
$ws.Range("K136").Value = 10155.24
$ws.Range("L136").Value = 11999.7276
$ws.Range("M136").Value = -7605.24
$ws.Range("N136").Value = -17099.7276

Write-Output "Applied Anima_Profits updates"
